$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the AE7:AF56 "test (MSE)" / "train (MSE)" hybrid-model columns
# that were previously blank (these drove the AVERAGE/STDEV.S summary
# formulas in row 58-59 to #DIV/0!).
$ws.Range("AE7").Value = 0.00086340456197939001
$ws.Range("AF7").Value = 0.00605089522702258707
$ws.Range("AE8").Value = 0.00084502880735631326
$ws.Range("AF8").Value = 0.00629160401161839127
$ws.Range("AE9").Value = 0.00077195715375899098
$ws.Range("AF9").Value = 0.00754655392799528143
$ws.Range("AE10").Value = 0.00071656105194302497
$ws.Range("AF10").Value = 0.00858267965502999136
$ws.Range("AE11").Value = 0.0008258253827837721
$ws.Range("AF11").Value = 0.00718340178230410763
$ws.Range("AE12").Value = 0.00073775063817157192
$ws.Range("AF12").Value = 0.00812832880889312499
$ws.Range("AE13").Value = 0.00080877247988693827
$ws.Range("AF13").Value = 0.00695132803993875076
$ws.Range("AE14").Value = 0.00087919170260952035
$ws.Range("AF14").Value = 0.00636687783892471264
$ws.Range("AE15").Value = 0.00089660909912945273
$ws.Range("AF15").Value = 0.00578424255601250794
$ws.Range("AE16").Value = 0.00079600310500722534
$ws.Range("AF16").Value = 0.00686765359621528595
$ws.Range("AE17").Value = 0.00087669525616632961
$ws.Range("AF17").Value = 0.0047843640386389932
$ws.Range("AE18").Value = 0.00071118319664440844
$ws.Range("AF18").Value = 0.00909075667769240932
$ws.Range("AE19").Value = 0.00088649397116871383
$ws.Range("AF19").Value = 0.00486870112484750778
$ws.Range("AE20").Value = 0.00091672137200346645
$ws.Range("AF20").Value = 0.00647739811173159759
$ws.Range("AE21").Value = 0.00080614299716547143
$ws.Range("AF21").Value = 0.00623702599125896141
$ws.Range("AE22").Value = 0.00079080541135545805
$ws.Range("AF22").Value = 0.0065147986790045501
$ws.Range("AE23").Value = 0.00086468352780674372
$ws.Range("AF23").Value = 0.00593623766252539091
$ws.Range("AE24").Value = 0.00085599131157726275
$ws.Range("AF24").Value = 0.00585201560983410283
$ws.Range("AE25").Value = 0.00084265711138683326
$ws.Range("AF25").Value = 0.00549088046715308538
$ws.Range("AE26").Value = 0.00083578475464690725
$ws.Range("AF26").Value = 0.00611635663969026504
$ws.Range("AE27").Value = 0.00092388903199000028
$ws.Range("AF27").Value = 0.00492886374362415218
$ws.Range("AE28").Value = 0.00091465408536584828
$ws.Range("AF28").Value = 0.00587417274403705461
$ws.Range("AE29").Value = 0.00085927450294244446
$ws.Range("AF29").Value = 0.00669572957235731528
$ws.Range("AE30").Value = 0.00084046554404440391
$ws.Range("AF30").Value = 0.00608098246615718964
$ws.Range("AE31").Value = 0.00083500849125374389
$ws.Range("AF31").Value = 0.00591862527380860627
$ws.Range("AE32").Value = 0.00078665944380502993
$ws.Range("AF32").Value = 0.00898024933669920762
$ws.Range("AE33").Value = 0.00084625978788563079
$ws.Range("AF33").Value = 0.00663925158507469842
$ws.Range("AE34").Value = 0.00081680227774547228
$ws.Range("AF34").Value = 0.00639672350533265158
$ws.Range("AE35").Value = 0.00078342656647699784
$ws.Range("AF35").Value = 0.00790211732195895446
$ws.Range("AE36").Value = 0.00089574455591543778
$ws.Range("AF36").Value = 0.00517296440950419743
$ws.Range("AE37").Value = 0.00082280602136381863
$ws.Range("AF37").Value = 0.00616844862821955362
$ws.Range("AE38").Value = 0.00074071045141117504
$ws.Range("AF38").Value = 0.01092020314707155965
$ws.Range("AE39").Value = 0.0008275783140295305
$ws.Range("AF39").Value = 0.00726483458748365066
$ws.Range("AE40").Value = 0.00084922984534872023
$ws.Range("AF40").Value = 0.00480287381610446269
$ws.Range("AE41").Value = 0.00076550630249273984
$ws.Range("AF41").Value = 0.00874732838830046336
$ws.Range("AE42").Value = 0.00084673677432635093
$ws.Range("AF42").Value = 0.00884047399879285545
$ws.Range("AE43").Value = 0.00080062766421530097
$ws.Range("AF43").Value = 0.00710220608234218842
$ws.Range("AE44").Value = 0.00090191750943323617
$ws.Range("AF44").Value = 0.00528087353036750962
$ws.Range("AE45").Value = 0.00082095438099731143
$ws.Range("AF45").Value = 0.0059319828540863196
$ws.Range("AE46").Value = 0.00084893170328617694
$ws.Range("AF46").Value = 0.00688128450103646198
$ws.Range("AE47").Value = 0.00074947183422892686
$ws.Range("AF47").Value = 0.00664851663807604007
$ws.Range("AE48").Value = 0.00091886300042368888
$ws.Range("AF48").Value = 0.00567488623134271026
$ws.Range("AE49").Value = 0.00087664367090482158
$ws.Range("AF49").Value = 0.00515143721111894007
$ws.Range("AE50").Value = 0.00081896554455745951
$ws.Range("AF50").Value = 0.0064955108341278071
$ws.Range("AE51").Value = 0.00092376417485305227
$ws.Range("AF51").Value = 0.00383698919978825697
$ws.Range("AE52").Value = 0.00081643309895520577
$ws.Range("AF52").Value = 0.00679559373283878526
$ws.Range("AE53").Value = 0.00091906258262140707
$ws.Range("AF53").Value = 0.00538893654349571066
$ws.Range("AE54").Value = 0.00083007365733830734
$ws.Range("AF54").Value = 0.0068132953316198076
$ws.Range("AE55").Value = 0.00073253315661487184
$ws.Range("AF55").Value = 0.00803114345146024519
$ws.Range("AE56").Value = 0.00079339724802980139
$ws.Range("AF56").Value = 0.00622394857119561336

# Update the saved view state: move the selection to V9 (this also clears
# the old topLeftCell scroll position recorded for the previous selection).
$ws.Range("V9").Select()
